# update 23 24 dil
# Adds event_log rows for the "Sloučení" (merge) episode 23 and the
# "Rivalové" tribe events of episode 24, renumbers the day value on
# rows 73-75 from 52 to 53, adds the matching personal_statistics row
# for Day 56 / episode 24, and tidies up the trailing blank rows /
# selections that Excel leaves behind.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("event_log")
$ws2 = $wb.Worksheets.Item("personal_statistics")

# --- event_log: day 52 -> 53 for the last three rows of episode 22 ---
$ws1.Cells.Item(73, 1).Value = 53
$ws1.Cells.Item(74, 1).Value = 53
$ws1.Cells.Item(75, 1).Value = 53

# --- event_log: row 76 - Day 54 / Episode 23 - merge ---
$ws1.Cells.Item(76, 1).Value = 54
$ws1.Cells.Item(76, 2).Value = 23
$ws1.Cells.Item(76, 3).Value = "11:28"
$ws1.Cells.Item(76, 4).Value = "Sloučení"
$ws1.Cells.Item(76, 6).Value = "Rivalové"
$ws1.Cells.Item(76, 7).Value = "Žaneta, Johanka, Kulhy, Adam, Martin, Barbora, Kristián, Pepa, Andrea, Karolína, Tomáš"

# --- event_log: row 77 - Day 56 / Episode 24 - personal immunity ---
$ws1.Cells.Item(77, 1).Value = 56
$ws1.Cells.Item(77, 2).Value = 24
$ws1.Cells.Item(77, 3).Value = "13:48"
$ws1.Cells.Item(77, 4).Value = "Souboj o osobní imunitu"
$ws1.Cells.Item(77, 6).Value = "Rivalové"
$ws1.Cells.Item(77, 7).Value = "Tomáš"

# --- event_log: row 78 - Day 57 / Episode 24 - kmenova rada ---
$ws1.Cells.Item(78, 1).Value = 57
$ws1.Cells.Item(78, 2).Value = 24
$ws1.Cells.Item(78, 3).Value = "51:02"
$ws1.Cells.Item(78, 4).Value = "Kmenová rada"
$ws1.Cells.Item(78, 6).Value = "Rivalové"
$ws1.Cells.Item(78, 7).Value = "Žaneta"
$ws1.Cells.Item(78, 8).Value = "Johanka"

# --- event_log: row 79 - Day 57 / Episode 24 - duel ---
$ws1.Cells.Item(79, 1).Value = 57
$ws1.Cells.Item(79, 2).Value = 24
$ws1.Cells.Item(79, 3).Value = "1:20:25"
$ws1.Cells.Item(79, 4).Value = "Duel"
$ws1.Cells.Item(79, 6).Value = "Rivalové"
$ws1.Cells.Item(79, 7).Value = "Johanka"
$ws1.Cells.Item(79, 8).Value = "Žaneta"

# --- event_log: drop the now-unused trailing blank row 103 ---
$ws1.Cells.Item(103, 3).Clear()

# --- event_log: row 91 no longer carries a styled-but-empty I cell ---
$ws1.Cells.Item(91, 9).Clear()

# --- personal_statistics: new Day 56 / Episode 24 scoring row ---
$ws2.Cells.Item(13, 1).Value = "Den 56"
$ws2.Cells.Item(13, 2).Value = 24
$ws2.Cells.Item(13, 3).Value = "Souboj o osobní imunitu"
$ws2.Cells.Item(13, 4).Value = 1
$ws2.Cells.Item(13, 5).Value = 2
$ws2.Cells.Item(13, 9).Value = 4
$ws2.Cells.Item(13, 10).Value = 8
$ws2.Cells.Item(13, 12).Value = 8
$ws2.Cells.Item(13, 14).Value = 4
$ws2.Cells.Item(13, 15).Value = 2
$ws2.Cells.Item(13, 16).Value = 4
$ws2.Cells.Item(13, 18).Value = 8
$ws2.Cells.Item(13, 22).Value = 8
$ws2.Cells.Item(13, 27).Value = 4

# --- restore sane selections on both sheets ---
[void]$ws2.Range("A9").Select()
$ws1.Select()
[void]$ws1.Range("A1").Select()
